$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vars_meta_data")
$ws.Name = "vars_meta_data_discrete"

$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item("vars_meta_data_discrete"))
$newSheet.Name = "vars_meta_data_pump"

$newSheet.Range("A1").Value = "h1"
$newSheet.Range("B1").Value = "h2"
Write-Output "check1"
Write-Output $newSheet.Range("A1").Value
Write-Output $newSheet.Range("B1").Value

$ws.Range("A1:J1").Copy()
$newSheet.Range("A1:J1").PasteSpecial(-4122)  # xlPasteFormats
Write-Output "check2"
Write-Output $newSheet.Range("A1").Value
Write-Output $newSheet.Range("B1").Value
